$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1251.0723
$ws.Range("I15").Value = 1251.0723
$ws.Range("K15").Value = 3753.2169
$ws.Range("M15").Value = -3584.2169
$ws.Range("H33").Value = 555885.25
$ws.Range("I33").Value = 1111321.1
$ws.Range("J33").Value = 449.44446
$ws.Range("K33").Value = 1111321.1
$ws.Range("L33").Value = 449.44446
$ws.Range("M33").Value = -1111092.1
$ws.Range("N33").Value = -907.4444599999999
$ws.Range("H70").Value = 3080.0
$ws.Range("I70").Value = 5000.0
$ws.Range("J70").Value = 2600.0
$ws.Range("K70").Value = 15000.0
$ws.Range("L70").Value = 7800.0
$ws.Range("M70").Value = -14730.0
$ws.Range("N70").Value = -8340.0
$ws.Range("H73").Value = 3080.0
$ws.Range("I73").Value = 5000.0
$ws.Range("J73").Value = 2600.0
$ws.Range("K73").Value = 15000.0
$ws.Range("L73").Value = 7800.0
$ws.Range("M73").Value = -14064.0
$ws.Range("N73").Value = -9672.0
$ws.Range("H100").Value = 1700.9048
$ws.Range("I100").Value = 1183.5294
$ws.Range("J100").Value = 3899.75
$ws.Range("K100").Value = 1183.5294
$ws.Range("L100").Value = 3899.75
$ws.Range("M100").Value = -642.5293999999999
$ws.Range("N100").Value = -4981.75
$ws.Range("H125").Value = 18682.5
$ws.Range("I125").Value = 22576.666
$ws.Range("J125").Value = 7000.0
$ws.Range("K125").Value = 203189.994
$ws.Range("L125").Value = 63000.0
$ws.Range("M125").Value = -200729.994
$ws.Range("N125").Value = -67920.0
$ws.Range("H134").Value = 57510.547
$ws.Range("J134").Value = 57510.547
$ws.Range("L134").Value = 57510.547
$ws.Range("N134").Value = -67650.54699999999
$ws.Range("H135").Value = 991.6667
$ws.Range("J135").Value = 950.0
$ws.Range("L135").Value = 8550.0
$ws.Range("N135").Value = -13620.0
$ws.Range("H136").Value = 54100.5
$ws.Range("J136").Value = 59840.7
$ws.Range("L136").Value = 59840.7
$ws.Range("N136").Value = -70040.7
$ws.Range("H138").Value = 6680.6777
$ws.Range("I138").Value = 4471.4287
$ws.Range("J138").Value = 6978.077
$ws.Range("K138").Value = 13414.2861
$ws.Range("L138").Value = 20934.231
$ws.Range("M138").Value = -8274.286100000001
$ws.Range("N138").Value = -31214.231
$ws.Range("H139").Value = 83832.375
$ws.Range("I139").Value = 40709.0
$ws.Range("K139").Value = 40709.0
$ws.Range("M139").Value = -35569.0
$ws.Range("H140").Value = 90278.57
$ws.Range("J140").Value = 90278.57
$ws.Range("L140").Value = 90278.57
$ws.Range("N140").Value = -100638.57

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 68125.0
$ws.Range("J7").Value = 90000.0
$ws.Range("L7").Value = 90000.0
$ws.Range("N7").Value = -90228.0
$ws.Range("H32").Value = 17476.547
$ws.Range("I32").Value = 10829.612
$ws.Range("J32").Value = 39189.867
$ws.Range("K32").Value = 10829.612
$ws.Range("L32").Value = 39189.867
$ws.Range("M32").Value = -10542.612
$ws.Range("N32").Value = -39763.867
$ws.Range("H52").Value = 54492.5
$ws.Range("J52").Value = 54492.5
$ws.Range("L52").Value = 54492.5
$ws.Range("N52").Value = -55128.5
$ws.Range("H117").Value = 65077.668
$ws.Range("J117").Value = 65077.668
$ws.Range("L117").Value = 65077.668
$ws.Range("N117").Value = -74255.668
$ws.Range("H118").Value = 80172.0
$ws.Range("J118").Value = 80172.0
$ws.Range("L118").Value = 80172.0
$ws.Range("N118").Value = -83486.0
$ws.Range("H122").Value = 1955.1538
$ws.Range("I122").Value = 1937.9546
$ws.Range("J122").Value = 2049.75
$ws.Range("K122").Value = 5813.8638
$ws.Range("L122").Value = 6149.25
$ws.Range("M122").Value = -3363.8638
$ws.Range("N122").Value = -11049.25
$ws.Range("H132").Value = 3408.24
$ws.Range("I132").Value = 3037.8125
$ws.Range("J132").Value = 4066.7778
$ws.Range("K132").Value = 9113.4375
$ws.Range("L132").Value = 12200.3334
$ws.Range("M132").Value = -6583.4375
$ws.Range("N132").Value = -17260.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 60140.777
$ws.Range("J20").Value = 9047.333
$ws.Range("L20").Value = 9047.333
$ws.Range("N20").Value = -9541.333
$ws.Range("H52").Value = 77916.0
$ws.Range("J52").Value = 77916.0
$ws.Range("L52").Value = 77916.0
$ws.Range("N52").Value = -78442.0
$ws.Range("H115").Value = 84284.57
$ws.Range("J115").Value = 93598.2
$ws.Range("L115").Value = 93598.2
$ws.Range("N115").Value = -96732.2
$ws.Range("H119").Value = 82397.164
$ws.Range("J119").Value = 82397.164
$ws.Range("L119").Value = 82397.164
$ws.Range("N119").Value = -92073.164
$ws.Range("H121").Value = 77916.0
$ws.Range("J121").Value = 77916.0
$ws.Range("L121").Value = 77916.0
$ws.Range("N121").Value = -81410.0
$ws.Range("H132").Value = 97454.14
$ws.Range("J132").Value = 97454.14
$ws.Range("L132").Value = 97454.14
$ws.Range("N132").Value = -107574.14
$ws.Range("H138").Value = 90613.7
$ws.Range("J138").Value = 90613.7
$ws.Range("L138").Value = 90613.7
$ws.Range("N138").Value = -100893.7
$ws.Range("H140").Value = 69248.09
$ws.Range("J140").Value = 69248.09
$ws.Range("L140").Value = 69248.09
$ws.Range("N140").Value = -79608.09

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1432.4615
$ws.Range("I16").Value = 1066.7142
$ws.Range("K16").Value = 1066.7142
$ws.Range("M16").Value = -779.7141999999999
$ws.Range("H18").Value = 28671.5
$ws.Range("J18").Value = 28671.5
$ws.Range("L18").Value = 28671.5
$ws.Range("N18").Value = -29131.5
$ws.Range("H22").Value = 10000000.0
$ws.Range("I22").Value = 10000000.0
$ws.Range("K22").Value = 10000000.0
$ws.Range("M22").Value = -9999650.0
$ws.Range("H99").Value = 1490453.6
$ws.Range("J99").Value = 3127612.5
$ws.Range("L99").Value = 3127612.5
$ws.Range("N99").Value = -3130608.5
$ws.Range("H107").Value = 1338.4783
$ws.Range("I107").Value = 1284.9445
$ws.Range("J107").Value = 1531.2
$ws.Range("K107").Value = 1284.9445
$ws.Range("L107").Value = 1531.2
$ws.Range("M107").Value = 635.0554999999999
$ws.Range("N107").Value = -5371.2
$ws.Range("H113").Value = 1432.4615
$ws.Range("I113").Value = 1066.7142
$ws.Range("K113").Value = 1066.7142
$ws.Range("M113").Value = 1103.2858
$ws.Range("H118").Value = 86733.664
$ws.Range("J118").Value = 86733.664
$ws.Range("L118").Value = 86733.664
$ws.Range("N118").Value = -90047.664
$ws.Range("H126").Value = 1490453.6
$ws.Range("J126").Value = 3127612.5
$ws.Range("L126").Value = 9382837.5
$ws.Range("N126").Value = -9387777.5
$ws.Range("H132").Value = 6249.5
$ws.Range("I132").Value = 5999.6665
$ws.Range("K132").Value = 17998.9995
$ws.Range("M132").Value = -15468.9995
$ws.Range("H138").Value = 69317.82
$ws.Range("J138").Value = 69981.11
$ws.Range("L138").Value = 69981.11
$ws.Range("N138").Value = -80261.11

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 3496.6667
$ws.Range("I22").Value = 5125.0
$ws.Range("K22").Value = 15375.0
$ws.Range("M22").Value = -15206.0
$ws.Range("H27").Value = 3496.6667
$ws.Range("I27").Value = 5125.0
$ws.Range("K27").Value = 15375.0
$ws.Range("M27").Value = -15273.0
$ws.Range("H98").Value = 307.33334
$ws.Range("I98").Value = 322.5
$ws.Range("J98").Value = 277.0
$ws.Range("K98").Value = 967.5
$ws.Range("L98").Value = 831.0
$ws.Range("M98").Value = 530.5
$ws.Range("N98").Value = -3827.0
$ws.Range("H122").Value = 1444201.2
$ws.Range("I122").Value = 1080.0
$ws.Range("J122").Value = 5052004.5
$ws.Range("K122").Value = 9720.0
$ws.Range("L122").Value = 45468040.5
$ws.Range("M122").Value = -7270.0
$ws.Range("N122").Value = -45472940.5
$ws.Range("H129").Value = 2010.5385
$ws.Range("I129").Value = 1932.8334
$ws.Range("J129").Value = 2077.1428
$ws.Range("K129").Value = 5798.5002
$ws.Range("L129").Value = 6231.428400000001
$ws.Range("M129").Value = -798.5002000000004
$ws.Range("N129").Value = -16231.4284

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1818.3529
$ws.Range("I102").Value = 1619.5
$ws.Range("K102").Value = 1619.5
$ws.Range("M102").Value = 2.5
$ws.Range("H122").Value = 12391.5
$ws.Range("I122").Value = 13241.071
$ws.Range("J122").Value = 6444.5
$ws.Range("K122").Value = 39723.213
$ws.Range("L122").Value = 19333.5
$ws.Range("M122").Value = -37273.213
$ws.Range("N122").Value = -24233.5
$ws.Range("H135").Value = 88483.336
$ws.Range("J135").Value = 88483.336
$ws.Range("L135").Value = 88483.336
$ws.Range("N135").Value = -98623.336
$ws.Range("H140").Value = 65467.65
$ws.Range("J140").Value = 64496.43
$ws.Range("L140").Value = 64496.43
$ws.Range("N140").Value = -74856.43

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value = 24210.0
$ws.Range("I94").Value = 10300.0
$ws.Range("J94").Value = 31165.0
$ws.Range("K94").Value = 10300.0
$ws.Range("L94").Value = 31165.0
$ws.Range("M94").Value = -9624.0
$ws.Range("N94").Value = -32517.0
$ws.Range("H117").Value = 40828.332
$ws.Range("J117").Value = 40828.332
$ws.Range("L117").Value = 40828.332
$ws.Range("N117").Value = -50006.332
$ws.Range("H123").Value = 80943.45
$ws.Range("J123").Value = 80943.45
$ws.Range("L123").Value = 80943.45
$ws.Range("N123").Value = -90743.45

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 6666866.5
$ws.Range("I113").Value = 6666866.5
$ws.Range("K113").Value = 20000599.5
$ws.Range("M113").Value = -19998429.5
